$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, sorted descending by 1989 value, with Russian and Uzbek (both 0) removed
$data = @(
    @("English", 29.96087130708418),
    @("Japanese", 8.976237444474988),
    @("Spanish", 8.494040116956178),
    @("German", 7.084007171339617),
    @("Chinese", 5.110516838379162),
    @("Arabic", 4.862741668193532),
    @("Portuguese", 4.470741487035319),
    @("French", 4.370839436412433),
    @("Italian", 4.24777423477317),
    @("Malay-Indonesian", 2.264699821659717),
    @("Dutch", 1.854872238053428),
    @("Persian", 1.400947752108371),
    @("Turkish", 1.366594729551016),
    @("Korean", 1.105953740502007),
    @("Polish", 1.028977090228833),
    @("Thai", 0.8092564987098517),
    @("Urdu", 0.7736562595404762),
    @("Swedish", 0.6407146613547712),
    @("Bengali", 0.4018064217702284),
    @("Vietnamese", 0.2779695048522767)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the two now-unused trailing rows (previously Uzbek / Vietnamese leftovers)
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
